# TeamContributions.xlsx - "Updated team spreadsheet for code review 2"
#
# The sheet used to track WEEK 1 TASKS (free text) + a SCORE column, one
# row per team member in an arbitrary order. For code review 2 it becomes
# a simple WEEK 1 / WEEK 2 numeric score tracker, reordered by member.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (also renames the Table's column headers) ---------------
$ws.Range("C2").Value = "WEEK 1 "
$ws.Range("D2").Value = "WEEK 2"

# --- Row 3: Joseph Fanning -----------------------------------------------
$ws.Range("B3").Value = "Joseph Fanning - 40593072"
$ws.Range("C3").Value = 33
$ws.Range("D3").Value = 33

# --- Row 4: Dan Ferguson --------------------------------------------------
$ws.Range("B4").Value = "Dan Ferguson - 40534169"
$ws.Range("C4").Value = 33
$ws.Range("D4").Value = 33

# --- Row 5: Callum Hamilton ----------------------------------------------
$ws.Range("B5").Value = "Callum Hamilton - 40591758"
$ws.Range("C5").Value = 33
$ws.Range("D5").Value = 33

# --- Row 6: Dale Follows ---------------------------------------------------
$ws.Range("B6").Value = "Dale Follows - 40606982"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0

# --- Column C is now a narrow score column instead of a wide text column --
$ws.Columns.Item(3).ColumnWidth = 15

# --- Selection left where the author's cursor ended up -------------------
$ws.Range("G11").Select()
